$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# 1) Main script sheet (ChaylaGreyJourney): wording tweaks
# -----------------------------------------------------------------
$journey = $wb.Worksheets.Item("ChaylaGreyJourney")
$journey.Range("B4").Value = "let go with me 🔥"
$journey.Range("B5").Value = "one sec"
$journey.Range("B8").Value = "jesus"
$journey.Range("B10").Value = "wait one sec"
$journey.Range("B11").Value = "I can't control myself anymore"
$journey.Range("B20").Value = "well?"
$journey.Range("B22").Value = "wait one sec"

# -----------------------------------------------------------------
# 2) Drop the original "boosters" tab; it gets rebuilt from scratch
#    below so the sheetId allocation order matches the target
#    workbook (cumcontrol1=32, cumcontrol2=33, dickpic=34, boosters=35).
# -----------------------------------------------------------------
$boosters = $wb.Worksheets.Item("boosters")
$boosters.Delete()

# -----------------------------------------------------------------
# 3) Duplicate "dickpic" (still unedited at this point) to become the
#    new "dickpic" tab, then immediately free the "dickpic" name by
#    renaming the original sheet to its new identity "cumcontrol2".
# -----------------------------------------------------------------
$cumcontrol2 = $wb.Worksheets.Item("dickpic")
$cumcontrol2.Copy($null, $cumcontrol2)
$cumcontrol2.Name = "cumcontrol2"
$newDickpic = $wb.Worksheets.Item("dickpic (2)")
$newDickpic.Name = "dickpic"

# -----------------------------------------------------------------
# 4) Duplicate the fresh "dickpic" once more to rebuild "boosters"
#    (still the original wording at this point).
# -----------------------------------------------------------------
$newDickpic.Copy($null, $newDickpic)
$newBoosters = $wb.Worksheets.Item("dickpic (2)")
$newBoosters.Name = "boosters"

# Rewrite the freshly-created "boosters" sheet with the original
# booster lines (unchanged content, just a fresh sheet part).
$newBoosters.Range("A2").Value = "h8"
$newBoosters.Range("B2").Value = "I came to America for a better life but right now the only thing I want is you 😏"
$newBoosters.Range("C2").Value = "BOOSTER. Chayla personality — Brazilian in USA."

$newBoosters.Range("A3").Value = "h7"
$newBoosters.Range("B3").Value = "more..."
$newBoosters.Range("C3").Value = "BOOSTER. Ultra micro."

$newBoosters.Range("A4").Value = "h6"
$newBoosters.Range("B4").Value = "meu deus my whole body is on fire 💕"
$newBoosters.Range("C4").Value = "BOOSTER. Physical + Portuguese flavor."

$newBoosters.Range("A5").Value = "h5"
$newBoosters.Range("B5").Value = "I literally can't think straight right now 🔥"
$newBoosters.Range("C5").Value = "BOOSTER."

$newBoosters.Range("A6").Value = "h4"
$newBoosters.Range("B6").Value = "you have no idea what you're doing to me"
$newBoosters.Range("C6").Value = "BOOSTER."

$newBoosters.Range("A7").Value = "h3"
$newBoosters.Range("B7").Value = "don't stop"
$newBoosters.Range("C7").Value = "BOOSTER. Micro."

$newBoosters.Range("A8").Value = "h2"
$newBoosters.Range("B8").Value = "I'm so wet right now because of you 😏"
$newBoosters.Range("C8").Value = "BOOSTER. Ego."

$newBoosters.Range("A9").Value = "h1"
$newBoosters.Range("B9").Value = "fuckkk 🔥"
$newBoosters.Range("C9").Value = "MID-SEXTING BOOSTER."

# -----------------------------------------------------------------
# 5) Rewrite the ORIGINAL "dickpic" sheet content (now named
#    "cumcontrol2") with the new delay/sync/edge names + wording.
# -----------------------------------------------------------------
$cumcontrol2.Range("A2").Value = "delay2"
$cumcontrol2.Range("B2").Value = "edge for me... just a little more... this last one is everything 😏"
$cumcontrol2.Range("C2").Value = "DELAY variant."

$cumcontrol2.Range("A3").Value = "delay1"
$cumcontrol2.Range("B3").Value = "hold it... what I'm about to send is the best one and you'll want to last for it"
$cumcontrol2.Range("C3").Value = "DELAY. Send PPV."

$cumcontrol2.Range("A4").Value = "sync2"
$cumcontrol2.Range("B4").Value = "I want us to finish at the same time... this one will push you over"
$cumcontrol2.Range("C4").Value = "SYNC variant."

$cumcontrol2.Range("A5").Value = "sync1"
$cumcontrol2.Range("B5").Value = "okay you earned it babe... let's go together, open this 🔥"
$cumcontrol2.Range("C5").Value = "SYNC. Send PPV."

$cumcontrol2.Range("A6").Value = "edge2"
$cumcontrol2.Range("B6").Value = "if you finish without my permission I'll be annoyed"
$cumcontrol2.Range("C6").Value = "EDGE variant."

$cumcontrol2.Range("A7").Value = "edge1"
$cumcontrol2.Range("B7").Value = "slow down babe, I'm in control here 😏"
$cumcontrol2.Range("C7").Value = "CONTROL."

# -----------------------------------------------------------------
# 6) Rewrite "cumcontrol" text (row names stay the same) -> becomes
#    "cumcontrol1"
# -----------------------------------------------------------------
$cumcontrol = $wb.Worksheets.Item("cumcontrol")
$cumcontrol.Range("B2").Value = "trust me you want to edge just a little longer for this one"

$cumcontrol.Range("B3").Value = "you're not done until I say you are... open this"
$cumcontrol.Range("C3").Value = "DELAY. Send PPV."

$cumcontrol.Range("B4").Value = "I'm right there too, let's finish this... but you need to see this first"
$cumcontrol.Range("C4").Value = "SYNC variant. Send PPV."

$cumcontrol.Range("B5").Value = "now... right now, with me babe. open this 🔥"
$cumcontrol.Range("C5").Value = "SYNC. Send PPV."

$cumcontrol.Range("B6").Value = "not a chance... you're going to wait until I say so 😏"

$cumcontrol.Range("B7").Value = "I didn't say you could cum yet babe 🔥"
$cumcontrol.Range("C7").Value = "CONTROL."

$cumcontrol.Name = "cumcontrol1"
